# "open dml running in parallel"
# Flip the "Testable" column (B) from "n" to "y" for the test case rows,
# leaving the header row (B1), the untouched block (B87:B90) and the
# already-"y" row (B113) alone. Also move the active selection to B88,
# matching the saved view state in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

foreach ($row in 2..86) {
    $ws.Cells.Item($row, 2).Value = "y"
}

foreach ($row in 91..112) {
    $ws.Cells.Item($row, 2).Value = "y"
}

$ws.Activate()
$ws.Range("B88").Select()
